$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = 0.16666666666666666
$ws.Range("G3").Value = 3.0
$ws.Range("L3").Value = 0.0
$ws.Range("M3").Value = 0.0
$ws.Range("N5").Value = 0.08695652173913043
$ws.Range("O5").Value = 4.0
$ws.Range("F6").Value = 0.225
$ws.Range("G6").Value = 9.0
$ws.Range("L6").Value = 0.05
$ws.Range("M6").Value = 2.0
$ws.Range("J7").Value = 0.0
$ws.Range("K7").Value = 0.0
$ws.Range("D8").Value = 0.017857142857142856
$ws.Range("E8").Value = 1.0
$ws.Range("J8").Value = 0.07142857142857142
$ws.Range("K8").Value = 4.0
$ws.Range("L8").Value = 0.10714285714285714
$ws.Range("M8").Value = 6.0
$ws.Range("N8").Value = 0.14285714285714285
$ws.Range("O8").Value = 8.0
$ws.Range("F9").Value = 0.03333333333333333
$ws.Range("G9").Value = 1.0
$ws.Range("H9").Value = 0.06666666666666667
$ws.Range("I9").Value = 2.0
$ws.Range("L9").Value = 0.2
$ws.Range("M9").Value = 6.0
$ws.Range("L11").Value = 0.029411764705882353
$ws.Range("M11").Value = 1.0
$ws.Range("N11").Value = 0.08823529411764706
$ws.Range("O11").Value = 3.0
$ws.Range("D12").Value = 0.07894736842105263
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 0.10526315789473684
$ws.Range("G12").Value = 4.0
$ws.Range("H12").Value = 0.18421052631578946
$ws.Range("I12").Value = 7.0
$ws.Range("H13").Value = 0.14285714285714285
$ws.Range("I13").Value = 4.0
$ws.Range("L13").Value = 0.10714285714285714
$ws.Range("M13").Value = 3.0
$ws.Range("L14").Value = 0.043478260869565216
$ws.Range("M14").Value = 1.0
$ws.Range("C15").Value = 33.0
$ws.Range("D15").Value = 0.030303030303030304
$ws.Range("F15").Value = 0.06060606060606061
$ws.Range("H15").Value = 0.12121212121212122
$ws.Range("J15").Value = 0.030303030303030304
$ws.Range("L15").Value = 0.030303030303030304
$ws.Range("N15").Value = 0.09090909090909091
$ws.Range("D16").Value = 0.06451612903225806
$ws.Range("E16").Value = 2.0
$ws.Range("H16").Value = 0.22580645161290322
$ws.Range("I16").Value = 7.0
$ws.Range("D18").Value = 0.2777777777777778
$ws.Range("E18").Value = 5.0
$ws.Range("H18").Value = 0.4444444444444444
$ws.Range("I18").Value = 8.0
$ws.Range("D19").Value = 0.28888888888888886
$ws.Range("E19").Value = 13.0
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 15.0
$ws.Range("H19").Value = 0.4888888888888889
$ws.Range("I19").Value = 22.0
$ws.Range("N19").Value = 0.15555555555555556
$ws.Range("O19").Value = 7.0
$ws.Range("D20").Value = 0.22
$ws.Range("E20").Value = 11.0
$ws.Range("J20").Value = 0.04
$ws.Range("K20").Value = 2.0
$ws.Range("L20").Value = 0.08
$ws.Range("M20").Value = 4.0
$ws.Range("D21").Value = 0.125
$ws.Range("E21").Value = 5.0
$ws.Range("F21").Value = 0.225
$ws.Range("G21").Value = 9.0
$ws.Range("L21").Value = 0.075
$ws.Range("M21").Value = 3.0
$ws.Range("D22").Value = 0.09302325581395349
$ws.Range("E22").Value = 4.0
$ws.Range("J22").Value = 0.11627906976744186
$ws.Range("K22").Value = 5.0
$ws.Range("H23").Value = 0.044444444444444446
$ws.Range("I23").Value = 2.0
$ws.Range("J23").Value = 0.1111111111111111
$ws.Range("K23").Value = 5.0
$ws.Range("D26").Value = 0.3220338983050847
$ws.Range("E26").Value = 19.0
$ws.Range("N26").Value = 0.4406779661016949
$ws.Range("O26").Value = 26.0
$ws.Range("D28").Value = 0.043478260869565216
$ws.Range("E28").Value = 1.0
$ws.Range("F28").Value = 0.21739130434782608
$ws.Range("G28").Value = 5.0
$ws.Range("D29").Value = 0.2
$ws.Range("E29").Value = 7.0
$ws.Range("H29").Value = 0.2571428571428571
$ws.Range("I29").Value = 9.0
$ws.Range("D34").Value = 0.029411764705882353
$ws.Range("E34").Value = 1.0
$ws.Range("N34").Value = 0.058823529411764705
$ws.Range("O34").Value = 2.0
$ws.Range("F36").Value = 0.0
$ws.Range("G36").Value = 0.0
$ws.Range("N36").Value = 0.3
$ws.Range("O36").Value = 6.0
$ws.Range("N37").Value = 0.047619047619047616
$ws.Range("O37").Value = 1.0
$ws.Range("N40").Value = 0.13793103448275862
$ws.Range("O40").Value = 4.0
$ws.Range("F41").Value = 0.14285714285714285
$ws.Range("G41").Value = 6.0
$ws.Range("H41").Value = 0.21428571428571427
$ws.Range("I41").Value = 9.0
$ws.Range("J41").Value = 0.047619047619047616
$ws.Range("K41").Value = 2.0
$ws.Range("L41").Value = 0.09523809523809523
$ws.Range("M41").Value = 4.0
$ws.Range("N41").Value = 0.14285714285714285
$ws.Range("O41").Value = 6.0
$ws.Range("F42").Value = 0.057971014492753624
$ws.Range("G42").Value = 4.0
$ws.Range("F43").Value = 0.06593406593406594
$ws.Range("G43").Value = 6.0
$ws.Range("H43").Value = 0.10989010989010989
$ws.Range("I43").Value = 10.0
$ws.Range("N43").Value = 0.07692307692307693
$ws.Range("O43").Value = 7.0
$ws.Range("F44").Value = 0.041666666666666664
$ws.Range("G44").Value = 1.0
$ws.Range("F46").Value = 0.13636363636363635
$ws.Range("G46").Value = 3.0
$ws.Range("H46").Value = 0.13636363636363635
$ws.Range("I46").Value = 3.0
$ws.Range("J46").Value = 0.13636363636363635
$ws.Range("K46").Value = 3.0
$ws.Range("H47").Value = 0.1388888888888889
$ws.Range("I47").Value = 5.0
$ws.Range("L48").Value = 0.1836734693877551
$ws.Range("M48").Value = 9.0
$ws.Range("N48").Value = 0.2857142857142857
$ws.Range("O48").Value = 14.0
$ws.Range("H49").Value = 0.19642857142857142
$ws.Range("I49").Value = 11.0
$ws.Range("L49").Value = 0.03571428571428571
$ws.Range("M49").Value = 2.0
$ws.Range("J50").Value = 0.02127659574468085
$ws.Range("K50").Value = 1.0
$ws.Range("L50").Value = 0.02127659574468085
$ws.Range("M50").Value = 1.0
$ws.Range("N50").Value = 0.06382978723404255
$ws.Range("O50").Value = 3.0
$ws.Range("D51").Value = 0.023809523809523808
$ws.Range("E51").Value = 1.0
$ws.Range("F51").Value = 0.023809523809523808
$ws.Range("G51").Value = 1.0
$ws.Range("L51").Value = 0.09523809523809523
$ws.Range("M51").Value = 4.0
$ws.Range("F52").Value = 0.1
$ws.Range("G52").Value = 4.0
$ws.Range("J52").Value = 0.0
$ws.Range("K52").Value = 0.0
$ws.Range("D53").Value = 0.08108108108108109
$ws.Range("E53").Value = 3.0
$ws.Range("J53").Value = 0.0
$ws.Range("K53").Value = 0.0
$ws.Range("N53").Value = 0.05405405405405406
$ws.Range("O53").Value = 2.0
$ws.Range("D54").Value = 0.05263157894736842
$ws.Range("E54").Value = 2.0
$ws.Range("F55").Value = 0.2413793103448276
$ws.Range("G55").Value = 7.0
$ws.Range("N55").Value = 0.10344827586206896
$ws.Range("O55").Value = 3.0
